# ---------------------------------------------------------------------------
# Adds a "Player Info" sheet in front of the existing "ODI Batting" /
# "ODI Bowling" sheets, and replaces the MATCH_CARD_LINK column (a full
# howstat.com scorecard URL) with a MATCH_CODE column (just the numeric
# match code) on both of the pre-existing sheets.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "Player Info" sheet, inserted before the current first sheet.
# ---------------------------------------------------------------------------
$battingSheetForAnchor = $wb.Worksheets.Item("ODI Batting")
$playerSheet  = $wb.Worksheets.Add($battingSheetForAnchor)
$playerSheet.Name = "Player Info"

# NOTE: the handle captured before Add() above ends up pointing at whatever
# sheet now sits in that slot (i.e. the newly inserted sheet) once the
# insert shifts everything else along -- so re-resolve "ODI Batting" by name
# afterwards rather than reusing $battingSheetForAnchor.
$battingSheet = $wb.Worksheets.Item("ODI Batting")

$playerSheet.Range("A1").Value = "ID"
$playerSheet.Range("B1").Value = "NAME"
$playerSheet.Range("C1").Value = "BATTING_HAND"
$playerSheet.Range("D1").Value = "BOWL_STYLE"

$headerRng = $playerSheet.Range("A1:D1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108  # xlCenter
$headerRng.VerticalAlignment = -4160    # xlTop
$headerRng.Borders.LineStyle = 1

# Keep the id as text ("4557") rather than a number -- match the source data's
# treatment of ids as plain strings. Writing with a leading apostrophe forces
# text, then resetting the style keeps it free of the quote-prefix format.
$playerSheet.Range("A2").Value = "'4557"
$playerSheet.Range("A2").Style = "Normal"
$playerSheet.Range("B2").Value = "Deepak Lokandersingh Chahar"
$playerSheet.Range("C2").Value = "Right Handed"
$playerSheet.Range("D2").Value = "Right Arm Medium"

# ---------------------------------------------------------------------------
# 2) "ODI Batting" sheet: MATCH_CARD_LINK (col D) -> MATCH_CODE, URL -> code.
# ---------------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingUsed = $battingSheet.UsedRange
$battingLastRow = $battingUsed.Rows.Count

for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $link = $cell.Value2
    if ($link -match 'MatchCode=(\d+)') {
        $cell.Value = "'" + $matches[1]
        $cell.Style = "Normal"
    }
}

# ---------------------------------------------------------------------------
# 3) "ODI Bowling" sheet: MATCH_CARD_LINK (col B) -> MATCH_CODE, URL -> code.
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingUsed = $bowlingSheet.UsedRange
$bowlingLastRow = $bowlingUsed.Rows.Count

for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $link = $cell.Value2
    if ($link -match 'MatchCode=(\d+)') {
        $cell.Value = "'" + $matches[1]
        $cell.Style = "Normal"
    }
}

Write-Output "Player Info sheet added; MATCH_CARD_LINK -> MATCH_CODE on both sheets."
